# ---------------------------------------------------------------------------
# Applies to "Actividades/portada.docx":
#   1. Wrap the whole cover-page text (title .. "2023") in a bookmark named
#      "_Hlk152000351" (bookmarkStart right before the title run,
#      bookmarkEnd right after the "2023" run).
#   2. Regenerate the wp14:editId of the cover picture's drawing anchor
#      (cosmetic Word "edit session" id -> 383CBEE9).
#   3. Drop the explicit paragraph-mark run formatting (<w:rPr> inside
#      <w:pPr>) that used to sit on the "2023" paragraph.
#   4. Delete the stray trailing empty paragraph at the end of the body.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 0: locate the paragraphs we need to touch. We search by content
# instead of hard-coding indices so the script keeps working even if
# paragraphs shift around.
#   $picParaIdx  -> paragraph containing the floating picture (wp:anchor)
#   $yearParaIdx -> paragraph whose whole text is "2023"
#   $lastParaIdx -> final paragraph of the document body
# ---------------------------------------------------------------------
$picParaIdx = -1
$yearParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.XML(0) -like "*<wp:anchor*") {
        $picParaIdx = $i
    }
    if ($p.Range.Text.Trim() -eq "2023") {
        $yearParaIdx = $i
    }
}
$lastParaIdx = $d.Paragraphs.Count

# ---------------------------------------------------------------------
# Step 1: regenerate the drawing's wp14:editId.
# The Word object model has no direct property for this cosmetic id, so
# we rewrite the paragraph's OOXML via InsertXML, changing only that
# attribute and leaving everything else (incl. the image relationship)
# untouched.
# ---------------------------------------------------------------------
$picPara = $d.Paragraphs($picParaIdx)
$picRange = $picPara.Range
$picXml = $picRange.XML(0)
$picXml = $picXml.Replace('wp14:editId="66791CC0"', 'wp14:editId="383CBEE9"')

# Pull just the <w:p>...</w:p> fragment for this paragraph back out of the
# full package XML returned by XML(0) and rebuild a minimal package around
# it so InsertXML only replaces this paragraph.
$pStartTag = $picXml.IndexOf("<w:p ")
if ($pStartTag -lt 0) { $pStartTag = $picXml.IndexOf("<w:p>") }
$pEndTag = $picXml.IndexOf("</w:p>", $pStartTag) + "</w:p>".Length
$picParaFragment = $picXml.Substring($pStartTag, $pEndTag - $pStartTag)

$picPackage = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:body>' + $picParaFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$picRange.InsertXML($picPackage)

# ---------------------------------------------------------------------
# Step 2: strip the paragraph-mark run properties (<w:rPr> living directly
# under <w:pPr>) from the "2023" paragraph, keeping its own run (and the
# run's formatting) unchanged. Rebuilt the same way, via InsertXML.
# ---------------------------------------------------------------------
$yearPara = $d.Paragraphs($yearParaIdx)
$yearRange = $yearPara.Range
$yearPackage = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr><w:t>2023</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$yearRange.InsertXML($yearPackage)

# ---------------------------------------------------------------------
# Step 3: wrap title..2023 in bookmark "_Hlk152000351" (start before the
# very first run of the document, end right after the "2023" run but
# before its paragraph mark).
# ---------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$yearPara = $d.Paragraphs($yearParaIdx)
$bmStart = $firstPara.Range.Start
$bmEnd = $yearPara.Range.End - 1
$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_Hlk152000351", $bmRange)

# ---------------------------------------------------------------------
# Step 4: delete the trailing empty paragraph at the end of the body
# (merges it into the "2023" paragraph, which keeps its own mark).
# ---------------------------------------------------------------------
$yearPara = $d.Paragraphs($yearParaIdx)
$lastPara = $d.Paragraphs($lastParaIdx)
if ($lastParaIdx -gt $yearParaIdx) {
    $trailingRange = $d.Range($yearPara.Range.End - 1, $lastPara.Range.End)
    $trailingRange.Delete()
}
